# Automatische test-sync: 2025-06-22 21:55:50
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Append new row 50 to the Logs sheet ---
$logs.Range("A50").Value = "Beschadigd product ontvangen"
$logs.Range("B50").Value = "mailmind.test@zohomail.eu"
$logs.Range("C50").Value = "Het product dat ik heb ontvangen is beschadigd aangekomen."
$logs.Range("D50").Value = "Retour / Terugbetaling"
$logs.Range("F50").Value = "2025-06-22 21:55:42"
$logs.Range("G50").Value = "Nee"

# --- Extend conditional formatting ranges to include the new row ---
$dFcs = $logs.Range("D2:D49").FormatConditions
for ($i = 1; $i -le $dFcs.Count; $i++) {
    $dFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D50"))
}

$gFcs = $logs.Range("G2:G49").FormatConditions
for ($i = 1; $i -le $gFcs.Count; $i++) {
    $gFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G50"))
}

# --- Update the Dashboard summary: Retour / Terugbetaling now leads with 6 ---
$dash.Range("A4").Value = "Retour / Terugbetaling"
$dash.Range("B4").Value = 6
$dash.Range("A5").Value = "Offerte / Prijsaanvraag"
$dash.Range("B5").Value = 5
